$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Q_cool (B3) to reflect new variable/fixed cost connection to demand model.
# This cascades into T_cool (B5 = B3/B4) via the existing formula.
$ws.Range("B3").Value = 36949.920440039998

$wb.Save()
